$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Recuperacao 01" (REC) grades in column D for students who took
# the recovery exam.
$ws.Range("D3").Value = 4
$ws.Range("D4").Value = 6
$ws.Range("D5").Value = 10
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 10
$ws.Range("D11").Value = 2
$ws.Range("D14").Value = 4
$ws.Range("D15").Value = 6
$ws.Range("D18").Value = 0
$ws.Range("D19").Value = 8
$ws.Range("D20").Value = 4
$ws.Range("D21").Value = 2
$ws.Range("D23").Value = 8
$ws.Range("D24").Value = 6
$ws.Range("D28").Value = 6
$ws.Range("D31").Value = 4
$ws.Range("D32").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 10
$ws.Range("D36").Value = 2
$ws.Range("D37").Value = 6
$ws.Range("D38").Value = 2
$ws.Range("D39").Value = 10

# Match the cursor / selection position left behind in the saved workbook.
$ws.Range("C17").Select()
